$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trade data appended below the existing rows (row 5)
$ws.Range("A5").Value = 42636.593159722222
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 10044.040000000001
$ws.Range("D5").Value = 10016
$ws.Range("E5").Value = 313.07
$ws.Range("F5").Value = 314.81
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 0.56000000000000005
$ws.Range("I5").Value = $false

# Match the date style used by the other rows in column A/G
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("G5").NumberFormat = "m/d/yy h:mm"

# Column C best-fit width grows slightly now that it holds the wider 10044.04 value
$ws.Columns("C").ColumnWidth = 8.166666666666666

